# Updated to fit config format changes
# Insert two new rows (block_num header row stays at 1; rows 2-9 unchanged).
# New rows are inserted right before the old row 10 ("trial_feedback"),
# pushing everything after it down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:11").Insert()

# New row 10: correct_kb_resp
$ws.Range("A10").Value = "correct_kb_resp"
$ws.Range("B10").Value = "string"
$ws.Range("C10").Value = "a-z"
$ws.Range("D10").Value = 'set to "None" for default (if trial_kb_resp = 0)'

# New row 11: trial_feedback (bool flag)
$ws.Range("A11").Value = "trial_feedback"
$ws.Range("B11").Value = "bool"
$ws.Range("C11").Value = "N/A"
$ws.Range("D11").Value = "0 = no trial feedback, 1 = display trial feedback"

# The old row 10 (now row 12, "trial_feedback") is renamed to "trial_feedback_type"
# -- its Data Type / Units / Meaning values stay the same.
$ws.Range("A12").Value = "trial_feedback_type"

# Restore the selection shown in the workbook.
$ws.Range("D9:D10").Select()
